$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "06_05_2020--23_48_22 535"
$ws.Range("H4").Value = "06_05_2020--23_49_16 646"
$ws.Range("H6").Value = "06_05_2020--23_49_48 956"
